$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue 'D2' '62.240.17'
Set-TextValue 'E2' '  -2.45%  '
Set-TextValue 'D3' '2.515.32'
Set-TextValue 'E3' '  -3.83%  '
Set-TextValue 'E4' '  +0.39%  '
Set-TextValue 'D5' '554.04'
Set-TextValue 'E5' '  -3.81%  '
Set-TextValue 'D6' '148.99'
Set-TextValue 'E6' '  -4.89%  '
Set-TextValue 'E7' '  +0.35%  '
Set-TextValue 'D8' '0.605'
Set-TextValue 'E8' '  -3.01%  '
Set-TextValue 'D9' '2.526.06'
Set-TextValue 'E9' '  -3.34%  '
Set-TextValue 'D10' '0.110'
Set-TextValue 'E10' '  -7.72%  '
Set-TextValue 'B11' 'TRON'
Set-TextValue 'C11' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D11' '0.155'
Set-TextValue 'E11' '  -1.05%  '
Set-TextValue 'B12' 'Toncoin'
Set-TextValue 'C12' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue 'D12' '5.43'
Set-TextValue 'E12' '  -6.82%  '
Set-TextValue 'D13' '0.362'
Set-TextValue 'E13' '  -5.20%  '
Set-TextValue 'D14' '26.46'
Set-TextValue 'E14' '  -6.42%  '
Set-TextValue 'D15' '2.992.99'
Set-TextValue 'E15' '  -3.12%  '
Set-TextValue 'D16' '0.0000169'
Set-TextValue 'E16' '  -7.71%  '
Set-TextValue 'D17' '62.096.98'
Set-TextValue 'E17' '  -2.38%  '
Set-TextValue 'D18' '2.544.69'
Set-TextValue 'E18' '  -2.49%  '
Set-TextValue 'D19' '11.40'
Set-TextValue 'E19' '  -5.88%  '
Set-TextValue 'D20' '7.09'
Set-TextValue 'E20' '  -7.97%  '
Set-TextValue 'D21' '4.26'
Set-TextValue 'E21' '  -6.51%  '
Set-TextValue 'D22' '325.13'
Set-TextValue 'E22' '  -5.22%  '
Set-TextValue 'D23' '0.998'
Set-TextValue 'E23' '  -0.09%  '
Set-TextValue 'D24' '65.55'
Set-TextValue 'E24' '  -2.83%  '
Set-TextValue 'D25' '1.75'
Set-TextValue 'E25' '  -0.59%  '
Set-TextValue 'D26' '0.0000104'
Set-TextValue 'E26' '  -4.89%  '
Set-TextValue 'D27' '2.687.23'
Set-TextValue 'E27' '  -2.40%  '
Set-TextValue 'D28' '8.63'
Set-TextValue 'E28' '  -6.61%  '
Set-TextValue 'D29' '1.51'
Set-TextValue 'E29' '  -4.84%  '
Set-TextValue 'B30' 'Binance-PegBSC-USD'
Set-TextValue 'C30' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D30' '1.00'
Set-TextValue 'E30' '  +0.23%  '
Set-TextValue 'B31' 'Bittensor'
Set-TextValue 'C31' 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D31' '537.86'
Set-TextValue 'E31' '  -9.78%  '
Set-TextValue 'D32' '7.70'
Set-TextValue 'E32' '  -2.95%  '
Set-TextValue 'D33' '0.154'
Set-TextValue 'E33' '  -3.81%  '
Set-TextValue 'D34' '1.92'
Set-TextValue 'E34' '  -6.83%  '
Set-TextValue 'D35' '1.62'
Set-TextValue 'E35' '  -8.81%  '
Set-TextValue 'D36' '5.99'
Set-TextValue 'E36' '  -9.49%  '
Set-TextValue 'D37' '4.94'
Set-TextValue 'E37' '  -8.25%  '
Set-TextValue 'E38' '  +0.39%  '
Set-TextValue 'D39' '0.381'
Set-TextValue 'E39' '  -6.48%  '
Set-TextValue 'D40' '18.78'
Set-TextValue 'E40' '  -5.06%  '
Set-TextValue 'D41' '151.19'
Set-TextValue 'E41' '  -2.11%  '
Set-TextValue 'D42' '1.74'
Set-TextValue 'E42' '  -7.05%  '
Set-TextValue 'E43' '  +0.14%  '
Set-TextValue 'D44' '2.31'
Set-TextValue 'E44' '  -5.86%  '
Set-TextValue 'D45' '150.21'
Set-TextValue 'E45' '  -3.71%  '
Set-TextValue 'B46' 'Filecoin'
Set-TextValue 'C46' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D46' '3.67'
Set-TextValue 'E46' '  -6.46%  '
Set-TextValue 'B47' 'InjectiveProtocol'
Set-TextValue 'C47' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D47' '22.10'
Set-TextValue 'E47' '  -5.44%  '
Set-TextValue 'D48' '0.0552'
Set-TextValue 'E48' '  -6.90%  '
Set-TextValue 'D49' '0.607'
Set-TextValue 'E49' '  -3.41%  '
Set-TextValue 'D50' '0.0957'
Set-TextValue 'E50' '  -6.02%  '
Set-TextValue 'D51' '0.0236'
Set-TextValue 'E51' '  -4.85%  '
